$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.545.00'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.070.29'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.51'
$ws.Range('E5').Value = '  +3.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.98'
$ws.Range('E6').Value = '  +4.21%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.434'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.26'
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('E11').Value = '  +3.60%  '
$ws.Range('D12').Value = '3.595.75'
$ws.Range('E12').Value = '  +3.18%  '
$ws.Range('E13').Value = '  +3.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.54'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('D16').Value = '57.598.84'
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = '3.072.03'
$ws.Range('E17').Value = '  +3.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.09'
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.87'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.11'
$ws.Range('E20').Value = '  +2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '331.73'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.77'
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('E25').Value = '  +4.74%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '0.0₃0902'
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.36'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('E31').Value = '  +4.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.74'
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '154.76'
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.48'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.08'
$ws.Range('E35').Value = '  +7.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.92'
$ws.Range('E36').Value = '  +4.20%  '
$ws.Range('E37').Value = '  +3.13%  '
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('D39').Value = '3.107.28'
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.93'
$ws.Range('E40').Value = '  +4.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.77'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.655'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').Value = '2.264.08'
$ws.Range('E44').Value = '  +5.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0258'
$ws.Range('E45').Value = '  +10.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '20.88'
$ws.Range('E46').Value = '  +7.48%  '
$ws.Range('E47').Value = '  +1.82%  '
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.924'
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '262.42'
$ws.Range('E50').Value = '  +15.90%  '
$ws.Range('E51').Value = '  +5.74%  '
